$d = $word.ActiveDocument

# Paragraph 19 (1-based 20 in Paragraphs collection concept, but using Content.Find covers all):
# 'אלמנט מסוג קלט מסוג תאריך, על מנת לקלט את התאריך להזמנה' -> 'אלמנט מסוג קלט מסוג תאריך'
$d.Content.Find.Execute(", על מנת לקלט את התאריך להזמנה", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Paragraph 20:
# 'אלמנט מסוג קלט של שעה, על מנת לקלוט מה השעה שהוזמנה' -> 'אלמנט מסוג מסוג שעה'
$d.Content.Find.Execute("קלט של שעה, על מנת לקלוט מה השעה שהוזמנה", $false, $false, $false, $false, $false, $true, 1, $false, "מסוג שעה", 2)

# Paragraph 21:
# 'אלמנט מסוג link שמקשר אייקון לאתר, בכדי שיהיה אייקון לאתר כשנפתח' -> 'אלמנט מסוג link שמקשר אייקון לאתר'
$d.Content.Find.Execute(", בכדי שיהיה אייקון לאתר כשנפתח", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Paragraph 23:
# 'מודל flexbox לסידור האלמנטים באתר בצורה של עמודה ובמרכז האתר' -> 'מודל flexbox לסידור האלמנטים באתר'
$d.Content.Find.Execute(" בצורה של עמודה ובמרכז האתר", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Paragraph 27:
# 'השתמשנו ב-camelCase על מנת שנוכל לקרוא בצורה שתהיה נוחה לעין את שמות המשתנים' -> 'השתמשנו ב-camelCase'
$d.Content.Find.Execute(" על מנת שנוכל לקרוא בצורה שתהיה נוחה לעין את שמות המשתנים", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Paragraph 28:
# 'עבור משתנים לא קבועים, השתמשנו ב-let על מנת לשמור על סקופ המשתנים בהתאם לפונק'' -> 'עבור משתנים לא קבועים, השתמשנו ב-let.'
$d.Content.Find.Execute(" על מנת לשמור על סקופ המשתנים בהתאם לפונק'", $false, $false, $false, $false, $false, $true, 1, $false, ".", 2)

# Paragraph 29:
# 'אינדנטציה - על מנת לשמור על נוחות קריאה של הפונק'' -> 'אינדנטציה'
$d.Content.Find.Execute(" - על מנת לשמור על נוחות קריאה של הפונק'", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

Write-Output "done"
